$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 58: 四方坪站 (站点 index 4)
$ws.Range("A58").Value = 45959
$ws.Range("B58").Value = "四方坪站"
$ws.Range("C58").Value = 8701.2800000000007
$ws.Range("D58").Value = 7029.51
$ws.Range("E58").Value = 2963.49
$ws.Range("F58").Value = 381

# Row 59: 高岭站 (站点 index 5)
$ws.Range("A59").Value = 45959
$ws.Range("B59").Value = "高岭站"
$ws.Range("C59").Value = 3363.35
$ws.Range("D59").Value = 2689.11
$ws.Range("E59").Value = 933.96
$ws.Range("F59").Value = 134

$ws.Range("J54").Select() | Out-Null
